# Update the "Pais" (COVID country stats) sheet per the 9-May-2020 02:04 data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Title timestamp update (row 1) ---
$ws.Range("A1").Value = 'Datos actualizados a 9 de Mayo de 2020 a las 02:04'

# --- Rows whose data overtook a neighbour in the "Casos totales" ranking. ---
# The refreshed totals push Nigeria / Sri Lanka / Republica de Africa Central /
# Montserrat above the countries that used to sit just ahead of them, so every
# row in the affected block is rewritten with its new country + full stat line.
$rankShiftRows = @{
    63 = @('Nigeria', 3912, 386, 679, 3116, 4, 10, 117)
    64 = @('Luxemburgo', 3871, 12, 3526, 245, 16, 0, 100)
    65 = @('Afganistan', 3778, 215, 472, 3197, 7, 3, 109)
    102 = @('Sri Lanka', 835, 12, 240, 586, 1, 0, 9)
    103 = @('Guatemala', 832, 34, 90, 719, 5, 2, 23)
    153 = @('Republica de Africa Central', 143, 49, 10, 133, 0, 0, 0)
    154 = @('Guayana Francesa', 141, 3, 113, 27, 0, 0, 1)
    155 = @('Brunei', 141, 0, 132, 8, 2, 0, 1)
    156 = @('Camboya', 122, 0, 120, 2, 1, 0, 0)
    157 = @('Sudan del Sur', 120, 46, 2, 118, 0, 0, 0)
    158 = @('Bermudas', 118, 0, 61, 50, 4, 0, 7)
    159 = @('Trinidad yTobago', 116, 0, 103, 5, 0, 0, 8)
    160 = @('Nepal', 102, 1, 31, 71, 0, 0, 0)
    161 = @('Uganda', 101, 0, 55, 46, 0, 0, 0)
    162 = @('Aruba', 101, 0, 89, 9, 4, 0, 3)
    163 = @('Monaco', 95, 0, 82, 9, 1, 0, 4)
    205 = @('Montserrat', 11, 0, 7, 3, 1, 0, 1)
    206 = @('Seychelles', 11, 0, 8, 3, 0, 0, 0)
}
foreach ($r in $rankShiftRows.Keys) {
    $vals = $rankShiftRows[$r]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $vals[$c]
    }
}

# --- Straight numeric refreshes; country + row position unchanged. ---
$cellUpdates = @{
    "B4" = 1321122
    "C4" = 28499
    "D4" = 223163
    "E4" = 1019382
    "F4" = 16938
    "G4" = 1649
    "H4" = 78577
    "B11" = 145892
    "C11" = 10199
    "E11" = 80550
    "G11" = 804
    "H11" = 9992
    "D15" = 30406
    "E15" = 31353
    "D49" = 4413
    "E49" = 3391
    "D144" = 187
    "E144" = 0
    "B145" = 186
    "C145" = 3
    "E145" = 89
    "B146" = 177
    "C146" = 1
    "E146" = 104
}
foreach ($addr in $cellUpdates.Keys) {
    $ws.Range($addr).Value = $cellUpdates[$addr]
}
